# Applies the edits described by the diff to the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "Field Explanation" --------------------------------------------
$wsField = $wb.Worksheets.Item("Field Explanation")

# Text content updates (shared-string text changes)
$wsField.Range("C9").Value  = "Used as 'bulk storage' for the boolean variables [Ownable][Owned][ContextMission][StaffPayed]"
$wsField.Range("C11").Value = "Eg.: 1011"
$wsField.Range("C22").Value = "Unused when flag is already true"
$wsField.Range("C34").Value = "If player can't meet StaffSal * Staff, no income is calculated"

# Update the selected cell shown when the sheet is active
$wsField.Activate()
$wsField.Range("C15").Select()

# --- Sheet "Properties Table" ---------------------------------------------
$wsProps = $wb.Worksheets.Item("Properties Table")

# Row 2 (Grotti) - income related fields now populated
$wsProps.Range("G2").Value = 100
$wsProps.Range("H2").Value = 1500000
$wsProps.Range("I2").Value = 6500
$wsProps.Range("J2").Value = 1500

# Row 3 (Perseus) - Flags value corrected
$wsProps.Range("C3").Value = 1001

# Update the selected range shown when the sheet is active
$wsProps.Activate()
$wsProps.Range("M2:M11").Select()
